# Generate Report for Handoff
# Update the "Latest Handoff Datetime" column (D) for the
# "a8046064-c542-4630-89ac-5b1091a93de8.md" row (row 5) on both the
# zh-cn and de-de worksheets to reflect the new handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-08 05:45:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-08 05:45:19"
